# Update cryptos list: prices (Price) and hourly volume change (Volume(1h))
# for rows 2-51, plus a position swap of the WhiteBITCoin/Stacks rows (48/49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'96.391.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "'3.702.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'235.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.76%  "
$ws.Range("E6").Value = "  -1.49%  "
$ws.Range("D7").Value = "'650.36"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.43%  "
$ws.Range("D8").Value = "'0.426"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  -4.09%  "
$ws.Range("D11").Value = "'3.699.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "'0.0000307"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +16.90%  "
$ws.Range("D13").Value = "'44.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "'6.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.82%  "
$ws.Range("D16").Value = "'4.392.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "'96.202.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'8.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.18%  "
$ws.Range("D19").Value = "'3.698.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").Value = "'13.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").Value = "'18.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "'0.503"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.13%  "
$ws.Range("D23").Value = "'520.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").Value = "'3.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("D26").Value = "'6.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").Value = "'101.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "'13.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("D29").Value = "'0.173"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.16%  "
$ws.Range("D30").Value = "'3.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("D31").Value = "'12.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  +6.96%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("D36").Value = "'32.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.87%  "
$ws.Range("D37").Value = "'648.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.07%  "
$ws.Range("D38").Value = "'0.585"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").Value = "'8.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'6.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.43%  "
$ws.Range("E42").Value = "  +4.85%  "
$ws.Range("D43").Value = "'40.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.35%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "'0.958"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").Value = "'0.0449"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("D47").Value = "'0.428"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.22%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "'23.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'2.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("D50").Value = "'8.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("D51").Value = "'3.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.66%  "
